$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "44.061.75"
$ws.Range("E2").Value = "  -0.97%  "
$ws.Range("D3").Value = "2.223.53"
$ws.Range("E3").Value = "  -0.68%  "
$ws.Range("E4").Value = "  -1.69%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "298.27"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.67%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "90.36"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.27%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.558"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.24%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.56%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.494"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -5.11%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "33.24"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.17%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0777"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.14%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.98"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.14%  "
$ws.Range("E13").Value = "  -0.74%  "
$ws.Range("D14").Value = "2.563.02"
$ws.Range("D15").Value = "2.226.47"
$ws.Range("E15").Value = "  +1.84%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "13.38"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.31%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.777"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -6.85%  "
$ws.Range("D18").Value = "44.006.51"
$ws.Range("E18").Value = "  -0.38%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.12"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.58%  "
$ws.Range("D20").Value = "0.0₃0906"
$ws.Range("E20").Value = "  -4.81%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.98"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -5.24%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "64.05"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.25%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "234.99"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.89%  "
$ws.Range("E24").Value = "  -4.87%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.999"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.28%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.84"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -7.16%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.26"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.54%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "39.01"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.58%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.38"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.34%  "
$ws.Range("B30").Value = "Monero"
$ws.Range("C30").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "151.75"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.86%  "
$ws.Range("B31").Value = "EthereumClassic"
$ws.Range("C31").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "19.21"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.70%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.48"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -7.95%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0762"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.02%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.49"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -5.96%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.116"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.94%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.103"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -6.25%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.84"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -6.88%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.68"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -6.63%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0299"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.45%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.60"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.26%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.15"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -7.06%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "13.50"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -9.65%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.00"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.62%  "
$ws.Range("D44").Value = "1.796.01"
$ws.Range("E44").Value = "  +0.18%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.79"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +8.65%  "
$ws.Range("E46").Value = "  -2.84%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "67.86"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.47%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "94.47"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.36%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "73.16"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -7.16%  "
$ws.Range("B50").Value = "FraxShare"
$ws.Range("C50").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.76"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.02%  "
$ws.Range("B51").Value = "THORChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.61"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -5.86%  "
